$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 12628025
$ws.Range("I46").Value = 5100
$ws.Range("J46").Value = 16835666
$ws.Range("K46").Value = 15300
$ws.Range("L46").Value = 50506998
$ws.Range("M46").Value = -15181
$ws.Range("N46").Value = -50507236

$ws.Range("H60").Value = 12628025
$ws.Range("I60").Value = 5100
$ws.Range("J60").Value = 16835666
$ws.Range("K60").Value = 15300
$ws.Range("L60").Value = 50506998
$ws.Range("M60").Value = -14816
$ws.Range("N60").Value = -50507966

$ws.Range("H98").Value = 7968.3335
$ws.Range("I98").Value = 6952.5
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 6952.5
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -5454.5
$ws.Range("N98").Value = -12996

$ws.Range("H112").Value = 1385.5
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 1407.4
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 4222.200000000001
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -6438.200000000001

$ws.Range("H122").Value = 7968.3335
$ws.Range("I122").Value = 6952.5
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 20857.5
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -18407.5
$ws.Range("N122").Value = -34900

$ws.Range("H129").Value = 751.6667
$ws.Range("I129").Value = 466.66666
$ws.Range("J129").Value = 1036.6666
$ws.Range("K129").Value = 1399.99998
$ws.Range("L129").Value = 3109.9998
$ws.Range("M129").Value = 3600.00002
$ws.Range("N129").Value = -13109.9998

$ws.Range("H134").Value = 63938.625
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 63938.625
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 63938.625
$ws.Range("N134").Value = -74078.625

$ws.Range("H137").Value = 4270.067
$ws.Range("I137").Value = 2225.125
$ws.Range("J137").Value = 6607.143
$ws.Range("K137").Value = 6675.375
$ws.Range("L137").Value = 19821.429
$ws.Range("M137").Value = -4125.375
$ws.Range("N137").Value = -24921.429

$ws.Range("H138").Value = 3428.9773
$ws.Range("I138").Value = 1819
$ws.Range("J138").Value = 4180.3
$ws.Range("K138").Value = 5457
$ws.Range("L138").Value = 12540.9
$ws.Range("M138").Value = -317
$ws.Range("N138").Value = -22820.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H74").Value = 60565.047
$ws.Range("I74").Value = 65927.62
$ws.Range("J74").Value = 18430.572
$ws.Range("K74").Value = 65927.62
$ws.Range("L74").Value = 18430.572
$ws.Range("M74").Value = -65053.62
$ws.Range("N74").Value = -20178.572

$ws.Range("H77").Value = 60565.047
$ws.Range("I77").Value = 65927.62
$ws.Range("J77").Value = 18430.572
$ws.Range("K77").Value = 329638.1
$ws.Range("L77").Value = 92152.86
$ws.Range("M77").Value = -325270.1
$ws.Range("N77").Value = -100888.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 20107
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20107
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 20107
$ws.Range("N61").Value = -20733

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1828.6666
$ws.Range("I16").Value = 1808.2858
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 1808.2858
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -1521.2858

$ws.Range("H31").Value = 3364.2666
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3364.2666
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3364.2666
$ws.Range("N31").Value = -3954.2666
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 3364.2666
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3364.2666
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3364.2666
$ws.Range("N34").Value = -3768.2666
$ws.Range("M34").ClearContents()

$ws.Range("H60").Value = 19166.666
$ws.Range("I60").Value = 14000
$ws.Range("J60").Value = 29500
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 29500
$ws.Range("M60").Value = -13489
$ws.Range("N60").Value = -30522

$ws.Range("H86").Value = 2084.7778
$ws.Range("I86").Value = 1537.2
$ws.Range("J86").Value = 2769.25
$ws.Range("K86").Value = 1537.2
$ws.Range("L86").Value = 2769.25
$ws.Range("M86").Value = -414.2
$ws.Range("N86").Value = -5015.25

$ws.Range("H89").Value = 2084.7778
$ws.Range("I89").Value = 1537.2
$ws.Range("J89").Value = 2769.25
$ws.Range("K89").Value = 7686
$ws.Range("L89").Value = 13846.25
$ws.Range("M89").Value = -2070
$ws.Range("N89").Value = -25078.25

$ws.Range("H99").Value = 2116.3333
$ws.Range("I99").Value = 2090.7693
$ws.Range("J99").Value = 2182.8
$ws.Range("K99").Value = 2090.7693
$ws.Range("L99").Value = 2182.8
$ws.Range("M99").Value = -592.7692999999999
$ws.Range("N99").Value = -5178.8

$ws.Range("H113").Value = 1828.6666
$ws.Range("I113").Value = 1808.2858
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1808.2858
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 361.7141999999999

$ws.Range("H122").Value = 9646.291999999999
$ws.Range("I122").Value = 5231.636
$ws.Range("J122").Value = 13381.77
$ws.Range("K122").Value = 15694.908
$ws.Range("L122").Value = 40145.31
$ws.Range("M122").Value = -13244.908
$ws.Range("N122").Value = -45045.31

$ws.Range("H126").Value = 2116.3333
$ws.Range("I126").Value = 2090.7693
$ws.Range("J126").Value = 2182.8
$ws.Range("K126").Value = 6272.3079
$ws.Range("L126").Value = 6548.400000000001
$ws.Range("M126").Value = -3802.3079
$ws.Range("N126").Value = -11488.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3466.6667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3466.6667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 10400.0001
$ws.Range("N88").Value = -11256.0001

$ws.Range("H91").Value = 3466.6667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3466.6667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 10400.0001
$ws.Range("N91").Value = -13364.0001

$ws.Range("H94").Value = 3545.4546
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3545.4546
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 10636.3638
$ws.Range("N94").Value = -11988.3638
$ws.Range("M94").ClearContents()

$ws.Range("H96").Value = 4261.2856
$ws.Range("I96").Value = 1101
$ws.Range("J96").Value = 4788
$ws.Range("K96").Value = 3303
$ws.Range("L96").Value = 14364
$ws.Range("M96").Value = -1244
$ws.Range("N96").Value = -18482

$ws.Range("H97").Value = 6652.048
$ws.Range("I97").Value = 856.1429000000001
$ws.Range("J97").Value = 9550
$ws.Range("K97").Value = 2568.4287
$ws.Range("L97").Value = 28650
$ws.Range("M97").Value = -2072.4287
$ws.Range("N97").Value = -29642

$ws.Range("H117").Value = 2000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 6000
$ws.Range("N117").Value = -12884
$ws.Range("M117").ClearContents()

$ws.Range("H129").Value = 1422.4
$ws.Range("I129").Value = 1125.8
$ws.Range("J129").Value = 1644.85
$ws.Range("K129").Value = 3377.4
$ws.Range("L129").Value = 4934.549999999999
$ws.Range("M129").Value = 1622.6
$ws.Range("N129").Value = -14934.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10284.857
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 10398.8
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10398.8
$ws.Range("M18").Value = -9707
$ws.Range("N18").Value = -10984.8

$ws.Range("H132").Value = 87073.62
$ws.Range("I132").Value = 127493.25
$ws.Range("J132").Value = 22402.2
$ws.Range("K132").Value = 382479.75
$ws.Range("L132").Value = 67206.60000000001
$ws.Range("M132").Value = -379949.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3654.1667
$ws.Range("I7").Value = 3350
$ws.Range("J7").Value = 4034.375
$ws.Range("K7").Value = 3350
$ws.Range("L7").Value = 4034.375
$ws.Range("M7").Value = -3238
$ws.Range("N7").Value = -4258.375

$ws.Range("H20").Value = 12000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 12000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12452

$ws.Range("H126").Value = 3654.1667
$ws.Range("I126").Value = 3350
$ws.Range("J126").Value = 4034.375
$ws.Range("K126").Value = 10050
$ws.Range("L126").Value = 12103.125
$ws.Range("M126").Value = -7580
$ws.Range("N126").Value = -17043.125

$ws.Range("H132").Value = 4213.6665
$ws.Range("I132").Value = 3996.8
$ws.Range("J132").Value = 4484.75
$ws.Range("K132").Value = 11990.4
$ws.Range("L132").Value = 13454.25
$ws.Range("M132").Value = -9460.400000000001
$ws.Range("N132").Value = -18514.25
